$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-05-03 Saturday"; new = "2025-05-04 Sunday"},
    @{old = "717×2=1434"; new = "433×2=866"},
    @{old = "285×9=2565"; new = "904×9=8136"},
    @{old = "845×7=5915"; new = "154×9=1386"},
    @{old = "643×5=3215"; new = "875×5=4375"},
    @{old = "888×6=5328"; new = "885×8=7080"},
    @{old = "889×2=1778"; new = "636×2=1272"},
    @{old = "410×3=1230"; new = "391×5=1955"},
    @{old = "962×7=6734"; new = "499×7=3493"},
    @{old = "734×5=3670"; new = "926×2=1852"},
    @{old = "589×2=1178"; new = "533×3=1599"},
    @{old = "132×8=1056"; new = "853×3=2559"},
    @{old = "450×6=2700"; new = "164×5=820"},
    @{old = "305×2=610"; new = "275×3=825"},
    @{old = "657×7=4599"; new = "440×5=2200"},
    @{old = "977×6=5862"; new = "684×7=4788"},
    @{old = "979×9=8811"; new = "611×3=1833"},
    @{old = "251×3=753"; new = "310×2=620"},
    @{old = "884×6=5304"; new = "367×3=1101"},
    @{old = "202×6=1212"; new = "831×6=4986"},
    @{old = "275×5=1375"; new = "818×8=6544"},
    @{old = "466×3=1398"; new = "434×7=3038"},
    @{old = "839×9=7551"; new = "741×6=4446"},
    @{old = "923×6=5538"; new = "531×4=2124"},
    @{old = "247×7=1729"; new = "328×6=1968"},
    @{old = "831×7=5817"; new = "388×8=3104"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
